# Apply automatic update of files: refresh the "Förändrad" date column (C)
# for all data rows, and re-sync the shuffled rows 6-16 with their
# (Beteckning / Datum / Area) triples as produced by the upstream data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 46072 to 46073 for every data row (2..16)
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46073
}

# Rows 6-16 had their (A, B, G) triples reassigned.
$rowData = @{
    6  = @{ A = "A 5792-2024";   B = 45335;               G = 5.6 }
    7  = @{ A = "A 8194-2025";   B = 45708;               G = 1.9 }
    8  = @{ A = "A 50997-2025";  B = 45946;               G = 1.5 }
    9  = @{ A = "A 12651-2022";  B = 44641;               G = 3.2 }
    10 = @{ A = "A 2593-2024";   B = 45313.69204861111;   G = 2.3 }
    11 = @{ A = "A 7333-2025";   B = 45703.35899305555;   G = 0.9 }
    13 = @{ A = "A 28288-2023";  B = 45099.6349537037;    G = 0.5 }
    15 = @{ A = "A 7827-2026";   B = 46062.63958333333;   G = 2.1 }
    16 = @{ A = "A 13651-2023";  B = 45006;               G = 2.2 }
}

foreach ($r in $rowData.Keys) {
    $data = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 7).Value = $data.G
}
